$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix capitalization of the stimulus filename: "stimuli/blank.JPG" -> "stimuli/blank.jpg"
# wherever it appears in the used range.
$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    $v = $cell.Value()
    if ($v -eq "stimuli/blank.JPG") {
        $cell.Value = "stimuli/blank.jpg"
    }
}

# Scroll the sheet so the top-left visible cell is C1 (was D1).
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1
